$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 10-14, column E (in_service) change from FALSE to TRUE
$ws.Range("E10:E14").Value = $true
